$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the affected columns keep their original Text storage so that
# values such as "245.02", "-0.57%", and numbers with leading/trailing
# zeros are written verbatim rather than being re-interpreted as
# numbers/percentages by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "245.02"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.57%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "28.50"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-4.53%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.229"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.08%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05695"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.52%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.610"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.43%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.194"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "3.12%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8503"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.66%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8565"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.47%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1369"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.26%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07034"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.41%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03134"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "6.57%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09205"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-1.92%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001535"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "1.66%"
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "One"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0005975"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.74%"
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005950"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-1.36%"
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.489"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.01%"
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.174"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-4.54%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.41%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.03279"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-4.87%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1287"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-2.27%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "0.29%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04087"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-1.77%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.1378"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.14%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001222"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.89%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-17.55%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-0.82%"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001447"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.33%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1064"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.87%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.003734"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-35.05%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002490"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "24.52%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009124"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-4.88%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005280"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "1.16%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000749"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.09%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1150"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "77.83%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-3.25%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002098"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.09%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0001998"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.09%"
